$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.046.78'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -1.47%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.667.33'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -1.28%  '

$ws.Range("E4").Value = '  +0.10%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.96'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.56%  '

$ws.Range("E6").Value = '  +0.37%  '

$ws.Range("E7").Value = '  +0.09%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2676'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +0.09%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06383'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +1.02%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.85'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -1.14%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07445'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.76%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.671.54'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -1.07%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.508'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -0.65%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5802'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +0.08%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.000008478'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -2.64%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.10'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -1.88%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '25.860.67'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -2.39%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '4.919'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -1.77%  '

$ws.Range("E19").Value = '  +0.03%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.79'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -1.14%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '189.21'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +1.11%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.187'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -1.45%  '

$ws.Range("E23").Value = '  +0.09%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '145.16'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.25%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '7.613'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +1.31%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1212'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +3.24%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.66'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -1.35%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.06653'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +15.34%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.329'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -1.57%  '

$ws.Range("E30").Value = '  -1.37%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.553'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +0.81%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.505'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -0.85%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.660'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -0.59%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.016'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +0.06%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.6162'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +3.12%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.369'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +0.52%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.685'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +0.24%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.323'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +7.16%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.097.27'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -0.38%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01592'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -1.39%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8679'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +0.66%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.009'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +0.71%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '101.47'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +1.40%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.815.71'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -1.57%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00000000116'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +3.06%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '56.22'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -0.64%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.005'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +0.02%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.089'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.90%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05231'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +0.20%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4283'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -0.59%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '5.985'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +2.56%  '
